$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1126.762
$ws.Range("J17").Value = 1126.762
$ws.Range("L17").Value = 3380.286
$ws.Range("N17").Value = -3716.286

$ws.Range("H40").Value = 3699.8
$ws.Range("I40").Value = 4899.6665
$ws.Range("K40").Value = 4899.6665
$ws.Range("M40").Value = -4724.6665

$ws.Range("H112").Value = 2263.9487
$ws.Range("I112").Value = 1000
$ws.Range("J112").Value = 2408.4
$ws.Range("K112").Value = 3000
$ws.Range("L112").Value = 7225.200000000001
$ws.Range("M112").Value = -1892
$ws.Range("N112").Value = -9441.200000000001

$ws.Range("H132").Value = 7941358.5
$ws.Range("I132").Value = 10422224
$ws.Range("K132").Value = 31266672
$ws.Range("M132").Value = -31264142

$ws.Range("H137").Value = 1082.7727
$ws.Range("I137").Value = 700.8461
$ws.Range("K137").Value = 2102.5383
$ws.Range("M137").Value = 447.4616999999998

$ws.Range("H138").Value = 479449.66
$ws.Range("I138").Value = 1323.4231
$ws.Range("J138").Value = 670700.1
$ws.Range("K138").Value = 3970.2693
$ws.Range("L138").Value = 2012100.3
$ws.Range("M138").Value = 1169.7307
$ws.Range("N138").Value = -2022380.3

$ws.Range("H141").Value = 566.96875
$ws.Range("I141").Value = 551.43335
$ws.Range("J141").Value = 800
$ws.Range("K141").Value = 1654.30005
$ws.Range("L141").Value = 2400
$ws.Range("M141").Value = 3525.69995
$ws.Range("N141").Value = -12760

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H133").Value = 38902.668
$ws.Range("J133").Value = 38902.668
$ws.Range("L133").Value = 38902.668
$ws.Range("N133").Value = -43962.668

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1643
$ws.Range("I20").Value = 1332.75
$ws.Range("K20").Value = 1332.75
$ws.Range("M20").Value = -1085.75

$ws.Range("H41").Value = 249800
$ws.Range("J41").Value = 249800
$ws.Range("L41").Value = 249800
$ws.Range("N41").Value = -250576

$ws.Range("H48").Value = 249800
$ws.Range("J48").Value = 249800
$ws.Range("L48").Value = 249800
$ws.Range("N48").Value = -250630

$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()

$ws.Range("H102").Value = 33999.8
$ws.Range("I102").Value = 18000
$ws.Range("J102").Value = 37999.75
$ws.Range("K102").Value = 18000
$ws.Range("L102").Value = 37999.75
$ws.Range("M102").Value = -14755
$ws.Range("N102").Value = -44489.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 175200
$ws.Range("I22").Value = 266.66666
$ws.Range("K22").Value = 266.66666
$ws.Range("M22").Value = 83.33334000000002

$ws.Range("H95").Value = 8000
$ws.Range("J95").Value = 8000
$ws.Range("L95").Value = 8000
$ws.Range("N95").Value = -13492

$ws.Range("H132").Value = 5738.4614
$ws.Range("I132").Value = 6754.0557
$ws.Range("J132").Value = 3453.375
$ws.Range("K132").Value = 20262.1671
$ws.Range("L132").Value = 10360.125
$ws.Range("M132").Value = -17732.1671
$ws.Range("N132").Value = -15420.125

$ws.Range("H134").Value = 1795.7391
$ws.Range("I134").Value = 1813.2632
$ws.Range("J134").Value = 1712.5
$ws.Range("K134").Value = 5439.7896
$ws.Range("L134").Value = 5137.5
$ws.Range("M134").Value = -2904.7896
$ws.Range("N134").Value = -10207.5

$ws.Range("H141").Value = 22000
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 9107.053
$ws.Range("I3").Value = 4210
$ws.Range("J3").Value = 15840.5
$ws.Range("K3").Value = 12630
$ws.Range("L3").Value = 47521.5
$ws.Range("M3").Value = -12518
$ws.Range("N3").Value = -47745.5

$ws.Range("H4").Value = 509116.2
$ws.Range("I4").Value = 224829.75
$ws.Range("K4").Value = 674489.25
$ws.Range("M4").Value = -674377.25

$ws.Range("H139").Value = 1914.6216
$ws.Range("J139").Value = 1699
$ws.Range("L139").Value = 5097
$ws.Range("N139").Value = -15377

$ws.Range("H140").Value = 32961.91
$ws.Range("I140").Value = 52458.1
$ws.Range("J140").Value = 2967.7693
$ws.Range("K140").Value = 157374.3
$ws.Range("L140").Value = 8903.3079
$ws.Range("M140").Value = -152194.3
$ws.Range("N140").Value = -19263.3079

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2589.658
$ws.Range("I102").Value = 2073.9666
$ws.Range("J102").Value = 4523.5
$ws.Range("K102").Value = 2073.9666
$ws.Range("L102").Value = 4523.5
$ws.Range("M102").Value = -451.9666000000002
$ws.Range("N102").Value = -7767.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 513.8095
$ws.Range("I16").Value = 513.8095
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 513.8095
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -343.8095
$ws.Range("N16").ClearContents()

$ws.Range("H22").Value = 1418.2
$ws.Range("I22").Value = 920.5
$ws.Range("J22").Value = 1750
$ws.Range("K22").Value = 920.5
$ws.Range("L22").Value = 1750
$ws.Range("M22").Value = -625.5
$ws.Range("N22").Value = -2340

$ws.Range("H27").Value = 1418.2
$ws.Range("I27").Value = 920.5
$ws.Range("J27").Value = 1750
$ws.Range("K27").Value = 920.5
$ws.Range("L27").Value = 1750
$ws.Range("M27").Value = -813.5
$ws.Range("N27").Value = -1964

$ws.Range("H46").Value = 1397.4286
$ws.Range("I46").Value = 697.5
$ws.Range("J46").Value = 2330.6667
$ws.Range("K46").Value = 697.5
$ws.Range("L46").Value = 2330.6667
$ws.Range("M46").Value = -509.5
$ws.Range("N46").Value = -2706.6667

$ws.Range("H47").Value = 6065
$ws.Range("J47").Value = 6065
$ws.Range("L47").Value = 6065
$ws.Range("N47").Value = -7045

$ws.Range("H52").Value = 6065
$ws.Range("J52").Value = 6065
$ws.Range("L52").Value = 6065
$ws.Range("N52").Value = -6531

$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()

$ws.Range("H88").Value = 11000
$ws.Range("I88").Value = 11000
$ws.Range("K88").Value = 11000
$ws.Range("M88").Value = -10572

$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()

$ws.Range("H91").Value = 11000
$ws.Range("I91").Value = 11000
$ws.Range("K91").Value = 11000
$ws.Range("M91").Value = -9518

$ws.Range("H132").Value = 29108.162
$ws.Range("I132").Value = 1644.3889
$ws.Range("J132").Value = 55126.473
$ws.Range("K132").Value = 4933.1667
$ws.Range("L132").Value = 165379.419
$ws.Range("M132").Value = -2403.1667
$ws.Range("N132").Value = -170439.419
